$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "2024-06-15 10:13:54"
$ws.Range("D26").Value = 200
$ws.Range("E26").Value = 7

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "2024-06-15 10:13:55"
$ws.Range("D27").Value = 200
$ws.Range("E27").Value = 2
